$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles, row structure) from column N into the new column O
# for rows 3-10 (row 1/2 header area spans only to N and is untouched).
$ws.Range("N3:N10").Copy($ws.Range("O3:O10"))

# Now set the new values for 2021 in column O.
$ws.Range("O4").Value = 2021
$ws.Range("O6").Value = 1860
$ws.Range("O7").Value = 1
$ws.Range("O8").Value = 510
$ws.Range("O9").Value = 178
$ws.Range("O10").Value = 821

# O3 and O5 remain empty (matching N3/N5 which only carry formatting).

# Update the active cell selection to match the authored workbook.
$ws.Range("P9").Select()
